$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("M2").Value = 493
$ws.Range("M3").Value = 541
$ws.Range("K4").Value = 1798
$ws.Range("M4").Value = 156
$ws.Range("M5").Value = 34
$ws.Range("M6").Value = 423
$ws.Range("K7").Value = 27592
$ws.Range("M7").Value = 1647

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("M3").Value = 39
$ws.Range("M7").Value = 112

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("M4").Value = 2
$ws.Range("M7").Value = 39

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("M4").Value = 4
$ws.Range("M5").Value = 4
$ws.Range("M6").Value = 20
$ws.Range("M7").Value = 58

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("M6").Value = 9
$ws.Range("M7").Value = 23

$ws = $wb.Worksheets.Item('New City')
$ws.Range("M3").Value = 14
$ws.Range("M7").Value = 37

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("M6").Value = 9
$ws.Range("M7").Value = 40

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("M3").Value = 1
$ws.Range("M7").Value = 5

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("M2").Value = 15
$ws.Range("M6").Value = 11
$ws.Range("M7").Value = 50
$ws.Range("M8").Value = 112
$ws.Range("M11").Value = 20
$ws.Range("M14").Value = 10
$ws.Range("M20").Value = 53
$ws.Range("M23").Value = 17
$ws.Range("M30").Value = 5
$ws.Range("L31").Value = 214
$ws.Range("M33").Value = 58
$ws.Range("M36").Value = 18
$ws.Range("M40").Value = 8
$ws.Range("M42").Value = 55
$ws.Range("M44").Value = 10
$ws.Range("M50").Value = 8
$ws.Range("L52").Value = 459
$ws.Range("M54").Value = 30
$ws.Range("K63").Value = 187
$ws.Range("M63").Value = 5
$ws.Range("M65").Value = 37
$ws.Range("M67").Value = 47
$ws.Range("M78").Value = 25
$ws.Range("M79").Value = 45
$ws.Range("M83").Value = 39
$ws.Range("M85").Value = 86
$ws.Range("M88").Value = 18
$ws.Range("M90").Value = 17
$ws.Range("M95").Value = 23
$ws.Range("M97").Value = 15
$ws.Range("M99").Value = 40
$ws.Range("K101").Value = 27592
$ws.Range("M101").Value = 1647

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("L4").Value = 12
$ws.Range("L7").Value = 214

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("M4").Value = 4
$ws.Range("M7").Value = 47

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("M6").Value = 14
$ws.Range("M7").Value = 30

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("M6").Value = 1
$ws.Range("M7").Value = 10

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("M2").Value = 3
$ws.Range("M6").ClearContents()
$ws.Range("M7").Value = 10

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("M3").Value = 6
$ws.Range("M7").Value = 11

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("M2").Value = 13
$ws.Range("M6").Value = 17
$ws.Range("M7").Value = 55

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("M3").Value = 12
$ws.Range("M7").Value = 25

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("M6").Value = 4
$ws.Range("M7").Value = 17

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("M3").Value = 16
$ws.Range("M6").Value = 10
$ws.Range("M7").Value = 45

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("M3").Value = 14
$ws.Range("M7").Value = 53

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("M3").Value = 7
$ws.Range("M7").Value = 18

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("M5").Value = 2
$ws.Range("M7").Value = 50

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("M4").Value = 3
$ws.Range("M7").Value = 8

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("M4").Value = 4
$ws.Range("M7").Value = 20

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("M6").Value = 5
$ws.Range("M7").Value = 15

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("M6").Value = 5
$ws.Range("M7").Value = 15

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("M2").Value = 5
$ws.Range("M7").Value = 18

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("M2").Value = 8
$ws.Range("M6").Value = 4
$ws.Range("M7").Value = 17

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("M2").Value = 25
$ws.Range("M3").Value = 38
$ws.Range("M4").Value = 4
$ws.Range("M7").Value = 86

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range("M3").Value = 2
$ws.Range("M6").Value = 3
$ws.Range("M7").Value = 8

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L2").Value = 146
$ws.Range("L7").Value = 459
